# Weekly refresh of fruit/vegetable (Hortaliza) price data.
# Each data row (2-24, row 7 unchanged) is updated in place with a
# "new" set of observations for columns D, H, J, K, L, M, N, O, P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44335; H = "Magnum";           J = 100; K = 35000; L = 36000; M = 35500; N = "`$/saco 25 kilos";  O = "Región Metropolitana";   P = 1420 }
    3  = @{ D = 44294; H = "Magnum";           J = 100; K = 24000; L = 25000; M = 24500; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 980  }
    4  = @{ D = 44384; H = "Sin especificar";  J = 100; K = 25000; L = 26000; M = 25500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1020 }
    5  = @{ D = 44435; H = "Magnum";           J = 100; K = 25000; L = 26000; M = 25500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1020 }
    6  = @{ D = 44188; H = "Magnum";           J = 100; K = 38000; L = 40000; M = 39000; N = "`$/saco 25 kilos";  O = "Región Metropolitana";   P = 1560 }
    8  = @{ D = 44272; H = "Magnum";           J = 100; K = 22000; L = 24000; M = 23000; N = "`$/saco 25 kilos";  O = "Región Metropolitana";   P = 920  }
    9  = @{ D = 44253; H = "Magnum";           J = 200; K = 25000; L = 26000; M = 25500; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 1020 }
    10 = @{ D = 44342; H = "Magnum";           J = 100; K = 28000; L = 30000; M = 29000; N = "`$/malla 25 kilos"; O = "Región Metropolitana";   P = 1160 }
    11 = @{ D = 44321; H = "Magnum";           J = 100; K = 24000; L = 25000; M = 24500; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 980  }
    12 = @{ D = 44399; H = "Magnum";           J = 100; K = 20000; L = 22000; M = 21000; N = "`$/malla 25 kilos"; O = "Perú";                   P = 840  }
    13 = @{ D = 44441; H = "Magnum";           J = 100; K = 28000; L = 29000; M = 28500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1140 }
    14 = @{ D = 44237; H = "Sin especificar";  J = 100; K = 20000; L = 22000; M = 21000; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 840  }
    15 = @{ D = 44167; H = "Sin especificar";  J = 100; K = 18000; L = 19000; M = 18500; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 740  }
    16 = @{ D = 44265; H = "Magnum";           J = 100; K = 20000; L = 22000; M = 21000; N = "`$/saco 25 kilos";  O = "Región Metropolitana";   P = 840  }
    17 = @{ D = 44447; H = "Magnum";           J = 100; K = 37000; L = 38000; M = 37500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1500 }
    18 = @{ D = 44433; H = "Magnum";           J = 100; K = 25000; L = 26000; M = 25500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1020 }
    19 = @{ D = 44279; H = "Magnum";           J = 100; K = 28000; L = 30000; M = 29000; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 1160 }
    20 = @{ D = 44363; H = "Magnum";           J = 100; K = 25000; L = 26000; M = 25500; N = "`$/malla 25 kilos"; O = "Perú";                   P = 1020 }
    21 = @{ D = 44160; H = "Magnum";           J = 100; K = 28000; L = 30000; M = 29000; N = "`$/malla 25 kilos"; O = "Región de O'Higgins";    P = 1160 }
    22 = @{ D = 44203; H = "Magnum";           J = 100; K = 20000; L = 22000; M = 21000; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 840  }
    23 = @{ D = 44244; H = "Magnum";           J = 100; K = 16000; L = 18000; M = 17000; N = "`$/saco 25 kilos";  O = "Región del Maule";       P = 680  }
    24 = @{ D = 44323; H = "Magnum";           J = 100; K = 20000; L = 22000; M = 21000; N = "`$/malla 25 kilos"; O = "Perú";                   P = 840  }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("D$r").Value = $row.D
    $ws.Range("H$r").Value = $row.H
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
}
